# Adapt column header formatting to respective input file names.
# - Rename the "_old" / "_new" suffixed headers to "_FV2404" / "_FV2410"
# - Freeze the header row (pane split after row 1)
# - Turn the data range into an Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row cells -------------------------------------------------
$oldSuffixHeaders = @(
    "Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID",
    "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung"
)

# Columns A..J -> "<Name>_FV2404"
for ($i = 0; $i -lt $oldSuffixHeaders.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value2 = "$($oldSuffixHeaders[$i])_FV2404"
}

# Column K stays "diff" (unchanged)

# Columns L..U -> "<Name>_FV2410"
for ($i = 0; $i -lt $oldSuffixHeaders.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value2 = "$($oldSuffixHeaders[$i])_FV2410"
}

# --- 2. Freeze the header row ----------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the data range into a table --------------------------------------
$range = $ws.Range("A1:U73")
$lo = $ws.ListObjects.Add(1, $range, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = $null | Out-Null
